# Commit: adding Following new scripts implementation,
#   1. User own profile comments test
#   2. Others Profile comments test
#
# Target worksheet: "Test Cases" (sheet1)
#   - Runmode column (C) for existing rows 2-7 changes from "Y" to "N"
#   - Two new test-case rows are appended (rows 8 and 9)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---------------------------------------------------------------------
# 1. Flip Runmode from Y to N for the existing test cases (rows 2-7)
# ---------------------------------------------------------------------
$ws.Range("C2:C7").Value = "N"

# ---------------------------------------------------------------------
# 2. Add the two new rows, copying the formatting of an existing plain
#    data row (row 5) so borders/number formats line up with the rest
#    of the table.
# ---------------------------------------------------------------------
$ws.Range("A5:D5").Copy($ws.Range("A8:D8"))
$ws.Range("A5:D5").Copy($ws.Range("A9:D9"))

$ws.Range("A8").Value = "OwnProfileCommentsLikeTest"
$ws.Range("B8").Value = "Validate User Own Profile Comments Like Test"
$ws.Range("C8").Value = "Y"
$ws.Range("D8").Value = "PASS"

$ws.Range("A9").Value = "OthersProfileCommentsLikeTest"
$ws.Range("B9").Value = "Validate Other User Profile Comments Like Test"
$ws.Range("C9").Value = "Y"
$ws.Range("D9").Value = "PASS"

# ---------------------------------------------------------------------
# 3. Widen the TCID column so the longer new test names fit (bestFit)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 34.140625

# ---------------------------------------------------------------------
# 4. Leave the selection where the author's last action landed
# ---------------------------------------------------------------------
[void]$ws.Range("C11").Select()
